$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 4 formatting/formulas down into row 5 (fill-down pattern used by the
# existing rows 2->3->4), then adjust the few cells that differ for the new
# case (wind speed 11 m/s, turbulence seed 101).
$ws.Range("A4:AM4").Copy($ws.Range("A5:AM5")) | Out-Null

# Relative formulas that follow on from row 4 (same pattern as rows 2-4)
$ws.Range("A5").Formula = "=A4"
$ws.Range("B5").Formula = "=A5+20"
$ws.Range("E5").Formula = "=E4+1"
$ws.Range("D5").Formula = "=""dlc01_steady_wsp"" & E5 & ""_s101"""
$ws.Range("I5").Formula = "=(0.16*(0.75*E5+5.6))/E5"
$ws.Range("K5").Formula = "=E5*B5/512"
$ws.Range("M5").Formula = "=8/E5"

# New turbulence seed name for the 11 m/s case
$ws.Range("J5").Value = "turb_s101_11ms"

# Narrow the sheet-tab / horizontal-scrollbar split in the window chrome
$excel.ActiveWindow.TabRatio = 500

$ws.Range("A5").Select()

$wb.Save()
